$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Update NO demand coefficients for the "Distributed Energy" / year 2040 row (row 4)
$ws.Range("R4").Value = 840
$ws.Range("S4").Value = 600
$ws.Range("T4").Value = 480
$ws.Range("U4").Value = 360
$ws.Range("V4").Value = 312

# Update the view: scroll and select cell Y4
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("Y4").Select()
